$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column C (Public Debt, GDP %) value updates, rows 2-22, plus new number format ---
$cValues = @{
    2  = 2430.7012405999999
    3  = 2496.1811048
    4  = 2640.2801350999998
    5  = 2759.6210000000001
    6  = 2740.2440000000001
    7  = 2731.971
    8  = 2766.7359999999999
    9  = 2895.2339999999999
    10 = 2983.2489999999998
    11 = 3029.8110000000001
    12 = 3138.8580000000002
    13 = 3289.2240000000002
    14 = 3503.4989999999998
    15 = 3751.2109999999998
    16 = 3958.0729999999999
    17 = 4058.5740000000001
    18 = 3867.8029999999999
    19 = 3950.607
    20 = 4028.489
    21 = 3999.0720000000001
    22 = 3970.7130000000002
}

foreach ($row in $cValues.Keys) {
    $cell = $ws.Range("C$row")
    $cell.Value = $cValues[$row]
    $cell.NumberFormat = "0.0000"
}

# Row 23: C23 gets a real value too (font differs slightly from above but reuse normal look)
$ws.Range("C23").Value = 4056.26

# --- Column Q (Interest Rate, %) updates ---
$qValues = @{
    4  = 11
    6  = 12.5
    7  = 14.9
    8  = 9.6
    9  = 5.43
    11 = 4.75
    12 = 2.5499999999999998
    13 = 2
    14 = 2.5
    15 = 2
    16 = 2.5
    17 = 3.5
    18 = 2.25
    19 = 1
    20 = 0.76
    21 = 0.76
    22 = 0.04
    23 = 0.04
}
foreach ($row in $qValues.Keys) {
    $ws.Range("Q$row").Value = $qValues[$row]
}
$ws.Range("Q5").Value = 11.3

# R23 gets a new value
$ws.Range("R23").Value = -2.1

# Row 8 gets a custom height
$ws.Rows.Item(8).RowHeight = 15.75

# Selection moves
$ws.Range("Q24").Select()
